$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone formatting from the row above (keeps the existing date/bool styles
# instead of minting new ones), then overwrite with the new trade's values.
$ws.Range("A3:I3").Copy($ws.Range("A4:I4"))

$ws.Range("A4").Value = 42633.676736111112
$ws.Range("B4").Value = $false
$ws.Range("C4").Value = 9948
$ws.Range("D4").Value = 10000
$ws.Range("E4").Value = 19.32
$ws.Range("F4").Value = 19.12
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = -1.04
$ws.Range("I4").Value = $false
